# Add a new "canonical SMILES" column (D) to the microstate table.
# For each molecule row, the new canonical SMILES is the existing
# "canonical isomeric SMILES" (column C) with the E/Z stereo-bond
# markers ('/' and '\\') stripped out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("D2").Value = "canonical SMILES"

# Data rows: 3 through 14
for ($r = 3; $r -le 14; $r++) {
    $smiles = $ws.Range("C$r").Value2
    $canonical = $smiles.Replace("/", "").Replace("\", "")
    $ws.Range("D$r").Value = $canonical
}

# Widen the new column to fit its contents
$ws.Columns.Item(4).ColumnWidth = 36
